$wb = $excel.ActiveWorkbook

# Video4 (BienBanHopVaLamViec sheet): no longer the active tab, selection moves to C3
$ws4 = $wb.Worksheets.Item("Video4")
$ws4.Activate()
$ws4.Range("C3").Select() | Out-Null

# Video5: add the new video link, becomes the active tab, selection moves to F5
$ws5 = $wb.Worksheets.Item("Video5")
$ws5.Range("A1").Value = "https://youtu.be/d8Ub1Z0KhAc"
$ws5.Activate()
$ws5.Range("F5").Select() | Out-Null
